$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 2 new rows for the extra "Periodo Mora" (2508) entries ------
# Before: rows 16-19 hold 4 data rows (2 workers x 2 periods: 2506, 2507).
# After:  rows 16-21 hold 6 data rows (2 workers x 3 periods: 2506, 2507, 2508).
# Inserting at 20:21 pushes the existing footer rows (24,25 -> 26,27) down
# together with their formatting/merged cells.
$ws.Rows("20:21").Insert()

# --- 2. Fix up borders/formatting for the data rows -------------------------
# Row 19 currently still carries the special "closing" border (thick bottom
# border) that used to mark the last row of the table. Move that formatting
# to the new last row (21), and give rows 19-20 the normal/open format that
# row 18 (an interior row) already has.
$ws.Range("B19:J19").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)
$ws.Range("B20:J20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 3. Rewrite the worker/period rows --------------------------------------
# New ordering: for each period (2506, 2507, 2508) list RICARDO then LUZ.
$ws.Range("C16").Value = "94430120"
$ws.Range("D16").Value = "RICARDO ANDRES OCAMPO ZABALA"
$ws.Range("E16").Value = "2506"

$ws.Range("C17").Value = "52184101"
$ws.Range("D17").Value = "LUZ YAMILE HERNANDEZ CARDENAS"
$ws.Range("E17").Value = "2506"

$ws.Range("C18").Value = "94430120"
$ws.Range("D18").Value = "RICARDO ANDRES OCAMPO ZABALA"
$ws.Range("E18").Value = "2507"

$ws.Range("C19").Value = "52184101"
$ws.Range("D19").Value = "LUZ YAMILE HERNANDEZ CARDENAS"
$ws.Range("E19").Value = "2507"

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "94430120"
$ws.Range("D20").Value = "RICARDO ANDRES OCAMPO ZABALA"
$ws.Range("E20").Value = "2508"
$ws.Range("F20").Value = 56940
$ws.Range("G20").Value = 1423500

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "52184101"
$ws.Range("D21").Value = "LUZ YAMILE HERNANDEZ CARDENAS"
$ws.Range("E21").Value = "2508"
$ws.Range("F21").Value = 56940
$ws.Range("G21").Value = 1423500

# --- 4. Update the summary figures ------------------------------------------
# Total "Valor Mora" now sums 6 rows of 56940 instead of 4.
$ws.Range("E11").Value = 341640
# One extra period added (2506, 2507, 2508 = 3 periods).
$ws.Range("F13").Value = 3
